# update new orleans xlsx files
#
# 1. Swap the tab order of the two sheets so "review_info" comes first and
#    "hotel_info" comes second.
# 2. Insert a new "State" column into "hotel_info" right after "Hotel_Name"
#    (before "City"), and populate the existing data row with "Louisiana".

$wb = $excel.ActiveWorkbook

# --- 1. Reorder sheets: review_info first, hotel_info second -------------
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($wb.Worksheets.Item(1))

# --- 2. Insert "State" column into hotel_info -----------------------------
$hotelSheet = $wb.Worksheets.Item("hotel_info")

# Hotel_Name is column B, City is column C -> insert a new column at C.
$hotelSheet.Columns.Item(3).Insert()

$hotelSheet.Cells.Item(1, 3).Value = "State"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"
